# fix the count issue
# - Remove the erroneous "祖籍..." clause and fix the "身後稱司馬溫公" wording
#   in 司馬光's biography (columns D and E of row 2).
# - Remove a spurious sentence about "毛斌公之子毛祥公" from 范純仁's
#   biography (column D of row 6).
# - Remove the stray "（yǐng）" pinyin annotation duplicated in 蘇轍's
#   biography (columns D and E of row 8).
# - Refresh the timestamp column (F) for rows 2-14 to reflect the new edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-SimaGuang([string]$text) {
    $text = $text.Replace("祖籍河内郡温县（今河南省焦作市温县），", "")
    $text = $text.Replace("身後稱司馬溫公", "因身後追封溫國公，又稱司馬溫公")
    return $text
}

# D2 / E2 - 司馬光
$d2 = $ws.Range("D2").Value2
$ws.Range("D2").Value = Fix-SimaGuang $d2

$e2 = $ws.Range("E2").Value2
$ws.Range("E2").Value = Fix-SimaGuang $e2

# D6 - 范純仁
$d6 = $ws.Range("D6").Value2
$ws.Range("D6").Value = $d6.Replace("毛斌公之子毛祥公是他的女婿。", "")

# D8 / E8 - 蘇轍
$d8 = $ws.Range("D8").Value2
$ws.Range("D8").Value = $d8.Replace("（yǐng）", "")

$e8 = $ws.Range("E8").Value2
$ws.Range("E8").Value = $e8.Replace("（yǐng）", "")

# Refresh timestamps in column F for rows 2-14
$ws.Range("F2").Value = "2024-09-28 02:04:46"
$ws.Range("F3").Value = "2024-09-28 02:04:47"
$ws.Range("F4").Value = "2024-09-28 02:04:48"
$ws.Range("F5").Value = "2024-09-28 02:04:49"
$ws.Range("F6").Value = "2024-09-28 02:04:50"
$ws.Range("F7").Value = "2024-09-28 02:04:51"
$ws.Range("F8").Value = "2024-09-28 02:04:52"
$ws.Range("F9").Value = "2024-09-28 02:04:53"
$ws.Range("F10").Value = "2024-09-28 02:04:54"
$ws.Range("F11").Value = "2024-09-28 02:04:55"
$ws.Range("F12").Value = "2024-09-28 02:04:56"
$ws.Range("F13").Value = "2024-09-28 02:04:57"
$ws.Range("F14").Value = "2024-09-28 02:04:57"
